$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")
$ws.Activate()

# Fix the mistranslated Spanish string in the last data row (D9):
# "selemente" -> "Buscar Producto" (correct translation of "Find Product").
$ws.Range("D9").Value = "Buscar Producto"

# Re-apply formatting to the last row (A9:D9), which produces a new cell
# style entry (same Arial font, default alignment/protection, just
# explicitly re-applied) distinct from the original default style.
$ws.Range("A9:D9").Style = "Normal"

# Move / restore the active selection to D9 (last populated cell).
$ws.Range("D9").Select()
